# Natmi following Dr Hou advice
# Rewrites the Nodal-Acvr2a LR-pairs table: the sending/target cluster set now
# includes "ECs" alongside "FAPs"/"sCs", expanding the data block from 3x2 (6 rows)
# to 3x3 (9 rows) and refreshing every numeric column (E:T) with recomputed stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nodal"
$ws.Range("C2").Value = "Acvr2a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3514346666666666
$ws.Range("H2").Value = 1.054304
$ws.Range("I2").Value = 0.164942595720082
$ws.Range("J2").Value = 0.1649425957200821
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.940402333333334
$ws.Range("N2").Value = 23.821207
$ws.Range("O2").Value = 0.1931648990487216
$ws.Range("P2").Value = 0.1931648990487216
$ws.Range("Q2").Value = 2.790532647214222
$ws.Range("R2").Value = 25.114793824928
$ws.Range("S2").Value = 0.03186111985110376
$ws.Range("T2").Value = 0.03186111985110376

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nodal"
$ws.Range("C3").Value = "Acvr2a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3514346666666666
$ws.Range("H3").Value = 1.054304
$ws.Range("I3").Value = 0.164942595720082
$ws.Range("J3").Value = 0.1649425957200821
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 26.95384733333333
$ws.Range("N3").Value = 80.861542
$ws.Range("O3").Value = 0.6557019380820612
$ws.Range("P3").Value = 0.6557019380820612
$ws.Range("Q3").Value = 9.472516352974221
$ws.Range("R3").Value = 85.25264717676799
$ws.Range("S3").Value = 0.1081531796859437
$ws.Range("T3").Value = 0.1081531796859437

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nodal"
$ws.Range("C4").Value = "Acvr2a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3514346666666666
$ws.Range("H4").Value = 1.054304
$ws.Range("I4").Value = 0.164942595720082
$ws.Range("J4").Value = 0.1649425957200821
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.212609666666666
$ws.Range("N4").Value = 18.637829
$ws.Range("O4").Value = 0.1511331628692172
$ws.Range("P4").Value = 0.1511331628692172
$ws.Range("Q4").Value = 2.183326407335111
$ws.Range("R4").Value = 19.649937666016
$ws.Range("S4").Value = 0.0249282961830346
$ws.Range("T4").Value = 0.02492829618303461

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nodal"
$ws.Range("C5").Value = "Acvr2a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.395965
$ws.Range("H5").Value = 4.187895
$ws.Range("I5").Value = 0.6551832032346963
$ws.Range("J5").Value = 0.6551832032346963
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.940402333333334
$ws.Range("N5").Value = 23.821207
$ws.Range("O5").Value = 0.1931648990487216
$ws.Range("P5").Value = 0.1931648990487216
$ws.Range("Q5").Value = 11.08452374325167
$ws.Range("R5").Value = 99.760713689265
$ws.Range("S5").Value = 0.1265583973112482
$ws.Range("T5").Value = 0.1265583973112482

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Nodal"
$ws.Range("C6").Value = "Acvr2a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.395965
$ws.Range("H6").Value = 4.187895
$ws.Range("I6").Value = 0.6551832032346963
$ws.Range("J6").Value = 0.6551832032346963
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 26.95384733333333
$ws.Range("N6").Value = 80.861542
$ws.Range("O6").Value = 0.6557019380820612
$ws.Range("P6").Value = 0.6557019380820612
$ws.Range("Q6").Value = 37.62662749267667
$ws.Range("R6").Value = 338.63964743409
$ws.Range("S6").Value = 0.4296048961598034
$ws.Range("T6").Value = 0.4296048961598034

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Nodal"
$ws.Range("C7").Value = "Acvr2a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.395965
$ws.Range("H7").Value = 4.187895
$ws.Range("I7").Value = 0.6551832032346963
$ws.Range("J7").Value = 0.6551832032346963
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.212609666666666
$ws.Range("N7").Value = 18.637829
$ws.Range("O7").Value = 0.1511331628692172
$ws.Range("P7").Value = 0.1511331628692172
$ws.Range("Q7").Value = 8.672585653328333
$ws.Range("R7").Value = 78.05327087995501
$ws.Range("S7").Value = 0.09901990976364476
$ws.Range("T7").Value = 0.09901990976364478

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Nodal"
$ws.Range("C8").Value = "Acvr2a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3832486666666666
$ws.Range("H8").Value = 1.149746
$ws.Range("I8").Value = 0.1798742010452218
$ws.Range("J8").Value = 0.1798742010452218
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.940402333333334
$ws.Range("N8").Value = 23.821207
$ws.Range("O8").Value = 0.1931648990487216
$ws.Range("P8").Value = 0.1931648990487216
$ws.Range("Q8").Value = 3.043148607046889
$ws.Range("R8").Value = 27.388337463422
$ws.Range("S8").Value = 0.03474538188636972
$ws.Range("T8").Value = 0.03474538188636972

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Nodal"
$ws.Range("C9").Value = "Acvr2a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3832486666666666
$ws.Range("H9").Value = 1.149746
$ws.Range("I9").Value = 0.1798742010452218
$ws.Range("J9").Value = 0.1798742010452218
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 26.95384733333333
$ws.Range("N9").Value = 80.861542
$ws.Range("O9").Value = 0.6557019380820612
$ws.Range("P9").Value = 0.6557019380820612
$ws.Range("Q9").Value = 10.33002605203689
$ws.Range("R9").Value = 92.970234468332
$ws.Range("S9").Value = 0.1179438622363142
$ws.Range("T9").Value = 0.1179438622363142

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Nodal"
$ws.Range("C10").Value = "Acvr2a"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3832486666666666
$ws.Range("H10").Value = 1.149746
$ws.Range("I10").Value = 0.1798742010452218
$ws.Range("J10").Value = 0.1798742010452218
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.212609666666666
$ws.Range("N10").Value = 18.637829
$ws.Range("O10").Value = 0.1511331628692172
$ws.Range("P10").Value = 0.1511331628692172
$ws.Range("Q10").Value = 2.380974371270444
$ws.Range("R10").Value = 21.428769341434
$ws.Range("S10").Value = 0.02718495692253781
$ws.Range("T10").Value = 0.02718495692253782
